$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-11-29 Friday"; New = "2024-11-30 Saturday" },
    @{ Old = "118÷9=13, 1";       New = "919÷6=153, 1" },
    @{ Old = "479÷7=68, 3";       New = "585÷9=65, 0" },
    @{ Old = "147÷8=18, 3";       New = "582÷2=291, 0" },
    @{ Old = "490÷3=163, 1";      New = "372÷4=93, 0" },
    @{ Old = "953÷4=238, 1";      New = "843÷4=210, 3" },
    @{ Old = "271÷9=30, 1";       New = "983÷4=245, 3" },
    @{ Old = "153÷7=21, 6";       New = "964÷6=160, 4" },
    @{ Old = "200÷7=28, 4";       New = "583÷3=194, 1" },
    @{ Old = "402÷5=80, 2";       New = "289÷2=144, 1" },
    @{ Old = "951÷8=118, 7";      New = "959÷3=319, 2" },
    @{ Old = "893÷4=223, 1";      New = "918÷8=114, 6" },
    @{ Old = "153÷6=25, 3";       New = "853÷7=121, 6" },
    @{ Old = "170÷3=56, 2";       New = "763÷2=381, 1" },
    @{ Old = "757÷3=252, 1";      New = "611÷4=152, 3" },
    @{ Old = "493÷5=98, 3";       New = "387÷4=96, 3" },
    @{ Old = "650÷8=81, 2";       New = "371÷9=41, 2" },
    @{ Old = "642÷5=128, 2";      New = "701÷6=116, 5" },
    @{ Old = "187÷8=23, 3";       New = "148÷3=49, 1" },
    @{ Old = "103÷4=25, 3";       New = "497÷3=165, 2" },
    @{ Old = "783÷6=130, 3";      New = "444÷4=111, 0" },
    @{ Old = "467÷4=116, 3";      New = "359÷5=71, 4" },
    @{ Old = "403÷3=134, 1";      New = "639÷4=159, 3" },
    @{ Old = "300÷2=150, 0";      New = "262÷9=29, 1" },
    @{ Old = "226÷9=25, 1";       New = "756÷3=252, 0" },
    @{ Old = "445÷2=222, 1";      New = "768÷4=192, 0" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2) | Out-Null
}

Write-Output "Done: applied $($replacements.Count) replacements"
